$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.576.17'
$ws.Range("E2").Value = '  +1.31%  '

# Row 3
$ws.Range("D3").Value = '1.767.89'
$ws.Range("E3").Value = '  -0.91%  '

# Row 4
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = '  +0.54%  '

# Row 5
$ws.Range("D5").Value = "'337.97"
$ws.Range("E5").Value = '  +0.93%  '

# Row 6
$ws.Range("D6").Value = "'1.007"
$ws.Range("E6").Value = '  +0.57%  '

# Row 7
$ws.Range("D7").Value = "'0.3858"
$ws.Range("E7").Value = '  +2.07%  '

# Row 8
$ws.Range("D8").Value = "'0.3417"
$ws.Range("E8").Value = '  -0.49%  '

# Row 9
$ws.Range("D9").Value = "'47.13"
$ws.Range("E9").Value = '  -2.51%  '

# Row 10
$ws.Range("D10").Value = "'1.147"
$ws.Range("E10").Value = '  -4.21%  '

# Row 11
$ws.Range("D11").Value = "'0.07447"
$ws.Range("E11").Value = '  -0.67%  '

# Row 12
$ws.Range("D12").Value = "'1.009"
$ws.Range("E12").Value = '  +0.72%  '

# Row 13
$ws.Range("D13").Value = "'22.60"
$ws.Range("E13").Value = '  +3.90%  '

# Row 14
$ws.Range("D14").Value = "'6.375"
$ws.Range("E14").Value = '  -1.52%  '

# Row 15
$ws.Range("D15").Value = '1.770.25'
$ws.Range("E15").Value = '  -0.93%  '

# Row 16
$ws.Range("D16").Value = "'7.074"
$ws.Range("E16").Value = '  -0.39%  '

# Row 17
$ws.Range("D17").Value = "'0.00001078"
$ws.Range("E17").Value = '  -1.80%  '

# Row 18
$ws.Range("D18").Value = "'0.06694"
$ws.Range("E18").Value = '  +0.36%  '

# Row 19
$ws.Range("D19").Value = "'82.48"
$ws.Range("E19").Value = '  -1.60%  '

# Row 20
$ws.Range("E20").Value = '  +0.59%  '

# Row 21
$ws.Range("D21").Value = "'17.44"
$ws.Range("E21").Value = '  +0.61%  '

# Row 22
$ws.Range("D22").Value = "'6.449"
$ws.Range("E22").Value = '  -2.55%  '

# Row 23
$ws.Range("D23").Value = '27.578.50'
$ws.Range("E23").Value = '  +1.33%  '

# Row 24
$ws.Range("D24").Value = "'12.18"
$ws.Range("E24").Value = '  -1.72%  '

# Row 25
$ws.Range("D25").Value = "'2.380"
$ws.Range("E25").Value = '  -1.52%  '

# Row 26
$ws.Range("D26").Value = "'20.82"
$ws.Range("E26").Value = '  -2.34%  '

# Row 27
$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").Value = "'1.431"
$ws.Range("E27").Value = '  -5.16%  '

# Row 28
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = "'2.442"
$ws.Range("E28").Value = '  -4.14%  '

# Row 29
$ws.Range("D29").Value = "'153.17"
$ws.Range("E29").Value = '  -0.48%  '

# Row 30
$ws.Range("B30").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C30").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D30").Value = '1.973.26'
$ws.Range("E30").Value = '  -0.81%  '

# Row 31
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = "'134.18"
$ws.Range("E31").Value = '  +0.14%  '

# Row 32
$ws.Range("D32").Value = "'6.157"
$ws.Range("E32").Value = '  +0.98%  '

# Row 33
$ws.Range("D33").Value = "'3.971"
$ws.Range("E33").Value = '  -1.18%  '

# Row 34
$ws.Range("D34").Value = "'0.08843"
$ws.Range("E34").Value = '  +1.67%  '

# Row 35
$ws.Range("D35").Value = "'12.68"
$ws.Range("E35").Value = '  -4.34%  '

# Row 36
$ws.Range("D36").Value = "'0.02438"
$ws.Range("E36").Value = '  +4.21%  '

# Row 37
$ws.Range("D37").Value = "'5.387"
$ws.Range("E37").Value = '  -1.30%  '

# Row 38
$ws.Range("D38").Value = "'0.6811"
$ws.Range("E38").Value = '  -2.04%  '

# Row 39
$ws.Range("D39").Value = "'0.06358"
$ws.Range("E39").Value = '  +0.35%  '

# Row 40
$ws.Range("D40").Value = "'0.2204"
$ws.Range("E40").Value = '  -0.01%  '

# Row 41
$ws.Range("D41").Value = "'1.544"
$ws.Range("E41").Value = '  -7.17%  '

# Row 42
$ws.Range("D42").Value = "'1.249"
$ws.Range("E42").Value = '  +0.46%  '

# Row 43
$ws.Range("D43").Value = "'8.410"
$ws.Range("E43").Value = '  -4.59%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = "'14.24"
$ws.Range("E44").Value = '  -0.86%  '

# Row 45
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").Value = "'1.007"
$ws.Range("E45").Value = '  +0.55%  '

# Row 46
$ws.Range("D46").Value = "'0.6263"
$ws.Range("E46").Value = '  -3.85%  '

# Row 47
$ws.Range("D47").Value = "'3.850"
$ws.Range("E47").Value = '  +0.03%  '

# Row 48
$ws.Range("D48").Value = "'132.09"
$ws.Range("E48").Value = '  +2.08%  '

# Row 49
$ws.Range("D49").Value = "'2.110"
$ws.Range("E49").Value = '  -1.77%  '

# Row 50
$ws.Range("D50").Value = "'0.07422"
$ws.Range("E50").Value = '  +4.06%  '

# Row 51
$ws.Range("D51").Value = "'1.240"
$ws.Range("E51").Value = '  +2.57%  '
